$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1080.753
$ws.Range("I15").Value = 1080.753
$ws.Range("K15").Value = 3242.259
$ws.Range("M15").Value = -3073.259

# Hunk 1: ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3257.1428
$ws.Range("I106").Value = 1950
$ws.Range("K106").Value = 1950
$ws.Range("M106").Value = -1319

# Hunk 2: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2797.697
$ws.Range("I137").Value = 1484.05
$ws.Range("J137").Value = 4818.6924
$ws.Range("K137").Value = 4452.15
$ws.Range("L137").Value = 14456.0772
$ws.Range("M137").Value = -1902.15
$ws.Range("N137").Value = -19556.0772

# Hunk 3: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4954.94
$ws.Range("I138").Value = 812.875
$ws.Range("J138").Value = 6904.147
$ws.Range("K138").Value = 2438.625
$ws.Range("L138").Value = 20712.441
$ws.Range("M138").Value = 2701.375
$ws.Range("N138").Value = -30992.441

# Hunk 4: ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6727.1113
$ws.Range("I141").Value = 6790.7427
$ws.Range("K141").Value = 20372.2281
$ws.Range("M141").Value = -15192.2281

# Hunk 5: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 870.0294
$ws.Range("I61").Value = 625.19354
$ws.Range("J61").Value = 3400
$ws.Range("K61").Value = 625.19354
$ws.Range("L61").Value = 3400
$ws.Range("M61").Value = -413.19354
$ws.Range("N61").Value = -3824

# Hunk 6: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3363.3142
$ws.Range("I74").Value = 3519.5557
$ws.Range("J74").Value = 2836
$ws.Range("K74").Value = 3519.5557
$ws.Range("L74").Value = 2836
$ws.Range("M74").Value = -2645.5557
$ws.Range("N74").Value = -4584

# Hunk 7: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3363.3142
$ws.Range("I77").Value = 3519.5557
$ws.Range("J77").Value = 2836
$ws.Range("K77").Value = 17597.7785
$ws.Range("L77").Value = 14180
$ws.Range("M77").Value = -13229.7785
$ws.Range("N77").Value = -22916

# Hunk 8: ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1984.5
$ws.Range("I122").Value = 1187.2667
$ws.Range("J122").Value = 3692.8572
$ws.Range("K122").Value = 3561.800099999999
$ws.Range("L122").Value = 11078.5716
$ws.Range("M122").Value = -1111.800099999999
$ws.Range("N122").Value = -15978.5716

# Hunk 9: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2232.7
$ws.Range("I132").Value = 1032.5714
$ws.Range("J132").Value = 5033
$ws.Range("K132").Value = 3097.7142
$ws.Range("L132").Value = 15099
$ws.Range("M132").Value = -567.7142000000003
$ws.Range("N132").Value = -20159

# Hunk 10: ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 26911.9
$ws.Range("J133").Value = 26911.9
$ws.Range("L133").Value = 26911.9
$ws.Range("N133").Value = -31971.9

# Hunk 11: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 870.0294
$ws.Range("I136").Value = 625.19354
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 1875.58062
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = 674.41938
$ws.Range("N136").Value = -15300

# Hunk 12: BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1044
$ws.Range("I107").Value = 994.63635
$ws.Range("J107").Value = 1111.875
$ws.Range("K107").Value = 994.63635
$ws.Range("L107").Value = 1111.875
$ws.Range("M107").Value = 925.36365
$ws.Range("N107").Value = -4951.875

# Hunk 13: CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4116.8125
$ws.Range("I132").Value = 3577.1
$ws.Range("J132").Value = 5016.3335
$ws.Range("K132").Value = 10731.3
$ws.Range("L132").Value = 15049.0005
$ws.Range("M132").Value = -8201.299999999999
$ws.Range("N132").Value = -20109.0005

# Hunk 14: CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7343.1055
$ws.Range("I134").Value = 8680.385
$ws.Range("J134").Value = 4445.6665
$ws.Range("K134").Value = 26041.155
$ws.Range("L134").Value = 13336.9995
$ws.Range("M134").Value = -23506.155
$ws.Range("N134").Value = -18406.9995

# Hunk 15: CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1559.2069
$ws.Range("I5").Value = 678.8889
$ws.Range("J5").Value = 1955.35
$ws.Range("K5").Value = 2036.6667
$ws.Range("L5").Value = 5866.049999999999
$ws.Range("M5").Value = -1924.6667
$ws.Range("N5").Value = -6090.049999999999

# Hunk 16: CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 25402.75
$ws.Range("J107").Value = 40389.96
$ws.Range("L107").Value = 121169.88
$ws.Range("N107").Value = -125009.88

# Hunk 17: CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 579.6053000000001
$ws.Range("I113").Value = 581.5417
$ws.Range("J113").Value = 576.2857
$ws.Range("K113").Value = 1744.6251
$ws.Range("L113").Value = 1728.8571
$ws.Range("M113").Value = 425.3749
$ws.Range("N113").Value = -6068.8571

# Hunk 18: CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2405.64
$ws.Range("I122").Value = 518.0833
$ws.Range("J122").Value = 3001.7104
$ws.Range("K122").Value = 4662.7497
$ws.Range("L122").Value = 27015.3936
$ws.Range("M122").Value = -2212.7497
$ws.Range("N122").Value = -31915.3936

# Hunk 19: CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10000867
$ws.Range("I131").Value = 71428870
$ws.Range("J131").Value = 958.6047
$ws.Range("K131").Value = 214286610
$ws.Range("L131").Value = 2875.8141
$ws.Range("M131").Value = -214281570
$ws.Range("N131").Value = -12955.8141

# Hunk 20: CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1559.2069
$ws.Range("I135").Value = 678.8889
$ws.Range("J135").Value = 1955.35
$ws.Range("K135").Value = 6110.0001
$ws.Range("L135").Value = 17598.15
$ws.Range("M135").Value = -3575.0001
$ws.Range("N135").Value = -22668.15

# Hunk 21: CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2889.1667
$ws.Range("I138").Value = 2395
$ws.Range("J138").Value = 3383.3333
$ws.Range("K138").Value = 7185
$ws.Range("L138").Value = 10149.9999
$ws.Range("M138").Value = -2045
$ws.Range("N138").Value = -20429.9999

# Hunk 22: CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2583.2917
$ws.Range("I139").Value = 1649.9286
$ws.Range("J139").Value = 3890
$ws.Range("K139").Value = 4949.7858
$ws.Range("L139").Value = 11670
$ws.Range("M139").Value = 190.2142000000003
$ws.Range("N139").Value = -21950

# Hunk 23: GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10996.5
$ws.Range("J5").Value = 10996.5
$ws.Range("L5").Value = 10996.5
$ws.Range("N5").Value = -11220.5

# Hunk 24: GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 725.1818
$ws.Range("I97").Value = 660
$ws.Range("K97").Value = 660
$ws.Range("M97").Value = -164

# Hunk 25: LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 93418.45
$ws.Range("I22").Value = 251600.25
$ws.Range("J22").Value = 3028.8572
$ws.Range("K22").Value = 251600.25
$ws.Range("L22").Value = 3028.8572
$ws.Range("M22").Value = -251305.25
$ws.Range("N22").Value = -3618.8572

# Hunk 26: LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 93418.45
$ws.Range("I27").Value = 251600.25
$ws.Range("J27").Value = 3028.8572
$ws.Range("K27").Value = 251600.25
$ws.Range("L27").Value = 3028.8572
$ws.Range("M27").Value = -251493.25
$ws.Range("N27").Value = -3242.8572

# Hunk 27: LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6242.08
$ws.Range("I40").Value = 4473.5454
$ws.Range("K40").Value = 4473.5454
$ws.Range("M40").Value = -4337.5454

# Hunk 28: LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1976.2222
$ws.Range("I136").Value = 1095.2778
$ws.Range("K136").Value = 3285.8334
$ws.Range("M136").Value = -735.8334000000004

# Hunk 29: WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 367334.34
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 550001.5
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 550001.5
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -551249.5

# Hunk 30: WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 367334.34
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 550001.5
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 2750007.5
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -2756247.5

# Hunk 31: WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2121.0625
$ws.Range("I126").Value = 1354.5217
$ws.Range("J126").Value = 4080
$ws.Range("K126").Value = 4063.5651
$ws.Range("L126").Value = 12240
$ws.Range("M126").Value = -1593.5651
$ws.Range("N126").Value = -17180

# Hunk 32: WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3444.4
$ws.Range("I136").Value = 1033.0667
$ws.Range("K136").Value = 3099.2001
$ws.Range("M136").Value = -549.2001
